# Commit: Minor fixes/scaling; Fixed race relations scaling by experience; Fewer encounters in space;

$wb = $excel.ActiveWorkbook

# --- 1. "Levels and Experience" sheet: move the selection to G15 -----------
$wsLE = $wb.Worksheets.Item("Levels and Experience")
$wsLE.Activate()
$wsLE.Range("G15").Select()

# --- 2. Add the new "Relations Levels" worksheet as the last tab ----------
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "Relations Levels"

# Explanatory text at the top of the sheet
$ws.Range("A1").Value = "Level = Math.Floor((Math.Sqrt(d + 0.25) - 0.5)"
$ws.Range("A2").Value = "Where d = exp / Scale"

# Header row (bold + underlined, matches existing "header" style used elsewhere)
# (column F is deliberately left untouched - it separates the two tables)
$ws.Range("B5:E5").Font.Bold = $true
$ws.Range("B5:E5").Font.Underline = $true
$ws.Range("G5:J5").Font.Bold = $true
$ws.Range("G5:J5").Font.Underline = $true
$ws.Range("B5").Value = "Exp"
$ws.Range("E5").Value = "Level"
$ws.Range("G5").Value = "Exp"
$ws.Range("J5").Value = "Level"

# Experience values that drive the two (positive / mirrored-negative) tables
$expValues = @(0,200,400,600,800,1000,1500,2000,2500,3000,3500,4000,4500,5000,6000,7000,8000,9000,10000,11000,12000,13000,14000,15000,16000,17000,18000,19000,20000,21000,22000,23000,24000,25000,26000,27000,28000,29000,30000,31000,32000,33000,34000,35000,40000)

$firstRow = 6
$lastRow = $firstRow + $expValues.Count - 1

$row = $firstRow
foreach ($v in $expValues) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# Row 6 gets plain (non-shared) formulas ...
$ws.Range("C6").Formula = "=(B6/1000) + 0.25"
$ws.Range("D6").Formula = "=SQRT(C6)-0.5"
$ws.Range("E6").Formula = "=FLOOR(D6,1)"
$ws.Range("G6").Formula = "=-B6"
$ws.Range("H6").Formula = "=(-G6/1000) + 0.25"
$ws.Range("I6").Formula = "=SQRT(H6)-0.5"
$ws.Range("J6").Formula = "=-FLOOR(I6,1)"

# ... rows 7-50 are filled as ranges so Excel collapses them into shared formulas
$rngC = "C7:C" + $lastRow
$rngD = "D7:D" + $lastRow
$rngE = "E7:E" + $lastRow
$rngG = "G7:G" + $lastRow
$rngH = "H7:H" + $lastRow
$rngI = "I7:I" + $lastRow
$rngJ = "J7:J" + $lastRow

$ws.Range($rngC).Formula = "=(B7/1000) + 0.25"
$ws.Range($rngD).Formula = "=SQRT(C7)-0.5"
$ws.Range($rngE).Formula = "=FLOOR(D7,1)"
$ws.Range($rngG).Formula = "=-B7"
$ws.Range($rngH).Formula = "=(-G7/1000) + 0.25"
$ws.Range($rngI).Formula = "=SQRT(H7)-0.5"
$ws.Range($rngJ).Formula = "=-FLOOR(I7,1)"

# Columns C:D and H:I are helper/working columns - hide them like the original
$ws.Range("C1:D1").EntireColumn.Hidden = $true
$ws.Range("H1:I1").EntireColumn.Hidden = $true

# Finally activate the new sheet and leave the given cell selected - this is
# also what makes it the workbook's active tab / tabSelected sheet.
$ws.Activate()
$ws.Range("N16").Select()
